$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 3 rows that correspond to the removed "ECs" sending-cluster pairs
# (original rows 8, 9, 10 after the earlier rows shift away) -- delete bottom-up
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

# Update remaining data rows (2-7) with refreshed TPM-derived values

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo1"
$ws.Range("C2").Value = "Lgr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.52558
$ws.Range("H2").Value = 4.57674
$ws.Range("I2").Value = 0.8891290081558957
$ws.Range("J2").Value = 0.8891290081558957
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.424719
$ws.Range("N2").Value = 4.274157
$ws.Range("O2").Value = 0.07423298812267187
$ws.Range("P2").Value = 0.07423298812267187
$ws.Range("Q2").Value = 2.17352281202
$ws.Range("R2").Value = 19.56170530818
$ws.Range("S2").Value = 0.06600270310195963
$ws.Range("T2").Value = 0.06600270310195963

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo1"
$ws.Range("C3").Value = "Lgr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.52558
$ws.Range("H3").Value = 4.57674
$ws.Range("I3").Value = 0.8891290081558957
$ws.Range("J3").Value = 0.8891290081558957
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.74065166666667
$ws.Range("N3").Value = 32.221955
$ws.Range("O3").Value = 0.5596266124066729
$ws.Range("P3").Value = 0.5596266124066729
$ws.Range("Q3").Value = 16.38572336963333
$ws.Range("R3").Value = 147.4715103267
$ws.Range("S3").Value = 0.4975802548267889
$ws.Range("T3").Value = 0.4975802548267889

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo1"
$ws.Range("C4").Value = "Lgr4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.52558
$ws.Range("H4").Value = 4.57674
$ws.Range("I4").Value = 0.8891290081558957
$ws.Range("J4").Value = 0.8891290081558957
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.027161333333335
$ws.Range("N4").Value = 21.081484
$ws.Range("O4").Value = 0.3661403994706553
$ws.Range("P4").Value = 0.3661403994706552
$ws.Range("Q4").Value = 10.72049678690667
$ws.Range("R4").Value = 96.48447108216001
$ws.Range("S4").Value = 0.3255460502271472
$ws.Range("T4").Value = 0.3255460502271471

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Rspo1"
$ws.Range("C5").Value = "Lgr4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.190234
$ws.Range("H5").Value = 0.570702
$ws.Range("I5").Value = 0.1108709918441043
$ws.Range("J5").Value = 0.1108709918441043
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.424719
$ws.Range("N5").Value = 4.274157
$ws.Range("O5").Value = 0.07423298812267187
$ws.Range("P5").Value = 0.07423298812267187
$ws.Range("Q5").Value = 0.271029994246
$ws.Range("R5").Value = 2.439269948214
$ws.Range("S5").Value = 0.008230285020712247
$ws.Range("T5").Value = 0.008230285020712247

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Rspo1"
$ws.Range("C6").Value = "Lgr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.190234
$ws.Range("H6").Value = 0.570702
$ws.Range("I6").Value = 0.1108709918441043
$ws.Range("J6").Value = 0.1108709918441043
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.74065166666667
$ws.Range("N6").Value = 32.221955
$ws.Range("O6").Value = 0.5596266124066729
$ws.Range("P6").Value = 0.5596266124066729
$ws.Range("Q6").Value = 2.043237129156667
$ws.Range("R6").Value = 18.38913416241
$ws.Range("S6").Value = 0.06204635757988396
$ws.Range("T6").Value = 0.06204635757988396

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Rspo1"
$ws.Range("C7").Value = "Lgr4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.190234
$ws.Range("H7").Value = 0.570702
$ws.Range("I7").Value = 0.1108709918441043
$ws.Range("J7").Value = 0.1108709918441043
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.027161333333335
$ws.Range("N7").Value = 21.081484
$ws.Range("O7").Value = 0.3661403994706553
$ws.Range("P7").Value = 0.3661403994706552
$ws.Range("Q7").Value = 1.336805009085334
$ws.Range("R7").Value = 12.031245081768
$ws.Range("S7").Value = 0.04059434924350813
$ws.Range("T7").Value = 0.04059434924350812
